$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 77: scene = kitchen2, Taito header
$ws.Range("A77").Value = "scene"
$ws.Range("B77").Value = "kitchen2"
$ws.Range("B77").NumberFormat = "@"
$ws.Range("E77").Value = "Taito"
$ws.Range("E77").Font.Bold = $true

# Row 78: algorithm = pt, N = 8
$ws.Range("A78").Value = "algorithm"
$ws.Range("B78").Value = "pt"
$ws.Range("B78").NumberFormat = "@"
$ws.Range("E78").Value = "N"
$ws.Range("F78").Value = 8

# Row 79: max depth = 17, n = 8
$ws.Range("A79").Value = "max depth"
$ws.Range("B79").Value = "17"
$ws.Range("B79").NumberFormat = "@"
$ws.Range("E79").Value = "n"
$ws.Range("F79").Value = 8

# Row 80: Box filter, cpus_per_task = 16
$ws.Range("A80").Value = "Box filter"
$ws.Range("E80").Value = "cpus_per_task"
$ws.Range("F80").Value = 16

# Row 81: cpus_total = F80*F79
$ws.Range("E81").Value = "cpus_total"
$ws.Range("F81").Formula = "=F80*F79"

# Row 82: time (m)
$ws.Range("E82").Value = "time (m)"
$ws.Range("F82").Value = 9.6265000000000001
$ws.Range("K82").Font.Bold = $true
$ws.Range("K82").Font.Bold = $false

# Row 83: time (s) = 60*F82
$ws.Range("E83").Value = "time (s)"
$ws.Range("F83").Formula = "=60*F82"

# Row 84: CPU eff. (%)
$ws.Range("E84").Value = "CPU eff. (%)"
$ws.Range("F84").Value = 75.48
$ws.Range("F84").Font.Bold = $true
$ws.Range("K84").Font.Bold = $true

# Row 85: Memory (MB)
$ws.Range("E85").Value = "Memory (MB)"
$ws.Range("F85").Value = 5350

# Row 87: SPP
$ws.Range("E87").Value = "SPP"
$ws.Range("F87").Value = "10k"

# Row 86: Network RCV
$ws.Range("E86").Value = "Network RCV"
$ws.Range("F86").Value = "12,2 M"

# Row 89: Note (bold), row 88 left blank
$ws.Range("E89").Value = "Bad util for simple pt"
$ws.Range("E89").Font.Bold = $true

# Update the view to match the target state
$excel.ActiveWindow.ScrollRow = 46
$ws.Range("I81").Select()
